$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D24").Value = "박효신 웃는 남자, 비를 피하는 나, 부동산은 나의 적"
$ws.Range("E24").Value = "https://blog.naver.com/hist0134/222819302056"

$ws.Range("D37").Value = "[Paper Review] The Rich Get Richer: Disparate Impact of Semi-Supervised Learning"
$ws.Range("E37").Value = "http://dsba.korea.ac.kr/seminar/?uid=2539&mod=document&pageid=1"

$ws.Range("D42").Value = "[Python] 코사인 유사도(Cosine Similarity)"
$ws.Range("E42").Value = "https://kjk92.tistory.com/92"

$ws.Range("D51").Value = "[pandas] 데이터프레임 컬럼 내 고유값의 개수 구하기, value_counts() 메소드"
$ws.Range("E51").Value = "https://bskyvision.com/entry/pandas-%EB%8D%B0%EC%9D%B4%ED%84%B0%ED%94%84%EB%A0%88%EC%9E%84-%EC%BB%AC%EB%9F%BC-%EB%82%B4-%EA%B3%A0%EC%9C%A0%EA%B0%92%EC%9D%98-%EA%B0%9C%EC%88%98-%EA%B5%AC%ED%95%98%EA%B8%B0-valuecounts-%EB%A9%94%EC%86%8C%EB%93%9C"
